# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values for rows 2-41 with the newly recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 4
    4  = 7
    5  = 7
    6  = 3
    7  = 8
    8  = 6
    9  = 8
    10 = 6
    11 = 7
    12 = 6
    13 = 3
    14 = 8
    15 = 8
    16 = 4
    17 = 4
    18 = 7
    19 = 8
    20 = 6
    21 = 7
    22 = 4
    23 = 8
    24 = 4
    25 = 5
    26 = 3
    27 = 9
    28 = 4
    29 = 2
    30 = 8
    31 = 3
    32 = 9
    33 = 6
    34 = 2
    35 = 6
    36 = 9
    37 = 3
    38 = 4
    39 = 2
    40 = 2
    41 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
